$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")

# Add new config row for db driver
$ws.Range("A8").Value = "dbdriver"
$ws.Range("B8").Value = "com.microsoft.sqlserver.jdbc.SQLServerDriver"

# Widen column B to fit the new longer value, keep C:D as before
$ws.Columns.Item(2).ColumnWidth = 41.29

# Update active selection
$ws.Range("B3").Select()

$wb.Save()
